$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 431, shifting existing rows 431..540 down to 432..541
$ws.Rows.Item(431).Insert()

# Populate the newly inserted row 431 with the new data record
$ws.Range("A431").Value = 9
$ws.Range("B431").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C431").Value = "Metropolitana"
$ws.Range("D431").Value = 44932
$ws.Range("E431").Value = 13
$ws.Range("F431").Value = 100112012
$ws.Range("G431").Value = "Espinaca"
$ws.Range("H431").Value = "Sin especificar"
$ws.Range("I431").Value = "Primera"
$ws.Range("J431").Value = 160
$ws.Range("K431").Value = 7000
$ws.Range("L431").Value = 8000
$ws.Range("M431").Value = 7500
$ws.Range("N431").Value = "$/cuna 10 kilos"
$ws.Range("O431").Value = "Provincia de Chacabuco"
$ws.Range("P431").Value = 750
$ws.Range("Q431").Value = 10
$ws.Range("R431").Value = "Hortaliza"
